$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the "Surface" time/accel block from D:E to H:I, and add new ---
# --- "Base"-side derived-quantity headers (accel/veloc/disp) in C:E    ---
# --- plus the matching derived-quantity headers for Surface in J:L.   ---

# Surface label moves from D1 to H1
$ws.Range("H1").Value = "Surface"
$ws.Range("D1").ClearContents()

# Surface time/accel headers move from D2:E2 to H2:I2
$ws.Range("H2").Value = "time (s)"
$ws.Range("I2").Value = "accel (g)"

$ws.Range("D2").ClearContents()
$ws.Range("D2").ClearFormats()
$ws.Range("E2").ClearContents()
$ws.Range("E2").ClearFormats()

# New derived-quantity headers for the "Base" block
$ws.Range("C2").Value = "accel (m/ss)"
$ws.Range("D2").Value = "v (m/s)"
$ws.Range("E2").Value = "d (m)"

# Matching derived-quantity headers for the "Surface" block
$ws.Range("J2").Value = "accel (m/ss)"
$ws.Range("K2").Value = "v (m/s)"
$ws.Range("L2").Value = "d (m)"

# --- Header formatting: bold everything, keep the per-column number formats ---

# time-like columns (A1/A2/H2) -> "0.00" + bold
foreach ($addr in @("A1", "A2", "H2")) {
    $ws.Range($addr).NumberFormat = "0.00"
    $ws.Range($addr).Font.Bold = $true
}

# accel-like columns (B1/B2/I2) -> "0.000000" + bold
foreach ($addr in @("B1", "B2", "I2")) {
    $ws.Range($addr).NumberFormat = "0.000000"
    $ws.Range($addr).Font.Bold = $true
}

# plain text headers -> General + bold
foreach ($addr in @("H1", "C2", "D2", "E2", "J2", "K2", "L2")) {
    $ws.Range($addr).Font.Bold = $true
}

# --- Column widths ---
$ws.Columns("A:A").ColumnWidth = 10.166666666666666
$ws.Columns("B:B").ColumnWidth = 10.166666666666666
$ws.Columns("C:C").ColumnWidth = 10.666666666666666
$ws.Columns("D:L").ColumnWidth = 10.166666666666666

# --- Page setup / view ---
$ws.PageSetup.Orientation = 1
$ws.Range("M2").Select()
